$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 : new bug - Bootstrap / special characters ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bootstrap"
$ws.Range("C5").Value = "L"

# Copy formatting from the row above (row 3) so the new cells pick up the
# same cell style used by the existing shared-formula rows.
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("D5").Formula = '=VLOOKUP(C5,$I$2:$J$5,2,FALSE)'
$ws.Range("E5").Formula = '=IF(G5 = "Y",0,D5)'

$ws.Range("F5").Value = "Bootstrap is unable to insert special characters"
$ws.Range("G5").Value = "N"

# --- Row 6 : new bug - UI - Links / broken links ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "UI - Links"
$ws.Range("C6").Value = "H"

$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("D6").Formula = '=VLOOKUP(C6,$I$2:$J$5,2,FALSE)'
$ws.Range("E6").Formula = '=IF(G6 = "Y",0,D6)'

$ws.Range("F6").Value = "Broken links"
$ws.Range("G6").Value = "N"

# Update selection to match the author's final cursor position
$ws.Range("G7").Select()
